$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has 8 "line" rows (rows 2-7, line1..line6) followed by
# 8 "extr" rows (rows 8-15, extr1..extr8). Two new rows, line7 and line8, are
# inserted right after line6 -- i.e. before extr1 -- which pushes extr1..extr8
# down from rows 8-15 to rows 10-17.

# 1) Shift the extr1..extr8 data (currently rows 8-15) down to rows 10-17.
#    Walk bottom-up so we never overwrite a source row before it's been read.
for ($r = 15; $r -ge 8; $r--) {
    $destRow = $r + 2
    for ($c = 1; $c -le 5; $c++) {
        $srcCell = $ws.Cells.Item($r, $c)
        $ws.Cells.Item($destRow, $c).Value = $srcCell.Value()
    }
}

# 2) The sheet grew by two rows (15 -> 17); rows 16:17 need the same formatting
#    (border/bold/centered style) as the rest of the data rows. Copy it over
#    from row 15, which still holds that formatting untouched.
$styleSrc = $ws.Range($ws.Cells.Item(15, 1), $ws.Cells.Item(15, 5))
$styleSrc.Copy()
$ws.Range($ws.Cells.Item(16, 1), $ws.Cells.Item(17, 5)).PasteSpecial(-4122)

# 3) Populate the two new rows: line7 (row 8) and line8 (row 9).
$ws.Cells.Item(8, 2).Value = "line7"
$ws.Cells.Item(8, 3).Value = 14
$ws.Cells.Item(8, 4).Value = 11
$ws.Cells.Item(8, 5).Value = $true

$ws.Cells.Item(9, 2).Value = "line8"
$ws.Cells.Item(9, 3).Value = 16
$ws.Cells.Item(9, 4).Value = 9
$ws.Cells.Item(9, 5).Value = $true

# 4) Re-sequence the index column (A) for all data rows (2-17) as 0..15.
for ($r = 2; $r -le 17; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}

# 5) Two in_service flags changed for what are now rows 14 (extr5) and 15 (extr6).
$ws.Cells.Item(14, 5).Value = $false
$ws.Cells.Item(15, 5).Value = $true
